$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParagraphIndexAt($pos) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $i
        }
    }
    return -1
}


# ---------------------------------------------------------------------------
# 1) Remove the standalone "Meta description: ..." paragraph that used to sit
#    right under the H1 title.
# ---------------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $metaParaIndex = Get-ParagraphIndexAt($findRange.Start)
    $metaPara = $d.Paragraphs.Item($metaParaIndex)
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Turn the closing italic "Create a feature image..." image-prompt
#    paragraph into two paragraphs:
#      - a new bold paragraph repeating the page title
#      - the same italic paragraph, now holding the meta-description text
# ---------------------------------------------------------------------------
$findRange2 = $d.Content
$found2 = $findRange2.Find.Execute("Create a feature image", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $targetIndex = Get-ParagraphIndexAt($findRange2.Start)
    $targetPara = $d.Paragraphs.Item($targetIndex)

    # Insert a new empty paragraph right before the target one, then stamp its
    # full OOXML content (leading empty run + bold run) in one shot so no
    # formatting is inherited from the neighbouring italic paragraph.
    $targetPara.Range.InsertParagraphBefore() | Out-Null
    $newPara = $d.Paragraphs.Item($targetIndex)
    $boldXml = "<w:p $wns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fortunes of Sparta Free - Exciting Slot Game!</w:t></w:r></w:p>"
    $newPara.Range.InsertXML($boldXml) | Out-Null

    # Replace the text of the (now shifted down by one) original paragraph.
    $finalPara = $d.Paragraphs.Item($targetIndex + 1)
    $italicXml = "<w:p $wns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Play Fortunes of Sparta for free and enjoy a high RTP percentage, Spartan Streak feature, and impressive graphics, animations, and sound effects.</w:t></w:r></w:p>"
    $finalPara.Range.InsertXML($italicXml) | Out-Null
}
